# "Removed some not needed components"
# Row 31 (C31): remove R58 from the 1k-resistor reference list.
# Row 36 (C36): remove R56 from the 100k-resistor reference list.
# Downstream formulas (counts, costs, totals) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 31: drop "R58," from the rich-text list, keep the other runs'
# ---- explicit colors (R39/R59 green, R64 red, the rest explicit black)
# ---- so the resulting rich-text run layout matches the original shape.
$cell31 = $ws.Range("C31")
$full31 = $cell31.Value()
$removeAt = $full31.IndexOf("R58,") + 1
$cell31.Characters($removeAt, 4).Text = ""
$new31 = $cell31.Value()

function Set-RunColor($cell, $fullText, $substring, $color) {
    $start = $fullText.IndexOf($substring) + 1
    $len = $substring.Length
    $cell.Characters($start, $len).Font.Color = $color
}

Set-RunColor $cell31 $new31 "R39" 5287936
Set-RunColor $cell31 $new31 ",R40,R50,R51,R57," 0
Set-RunColor $cell31 $new31 "R59" 5287936
Set-RunColor $cell31 $new31 ",R62," 0
Set-RunColor $cell31 $new31 "R64" 255
Set-RunColor $cell31 $new31 ",R65" 0

# ---- Row 36: drop ",R56" (plain text, no rich-text runs involved).
$cell36 = $ws.Range("C36")
$full36 = $cell36.Value()
$removeAt36 = $full36.IndexOf(",R56") + 1
$cell36.Characters($removeAt36, 4).Text = ""

# ---- Leave the view pointed at the row that was edited, like the author did.
$cell36.Select()
